# Apply a cyclic rotation of data among rows 11, 12 and 13:
#   new row 11 <= old row 13
#   new row 12 <= old row 11
#   new row 13 <= old row 12
# Columns affected: A, B, E, F, G, H, Q, R, Z, AB

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "E", "F", "G", "H", "Q", "R", "Z", "AB")

# Capture current (pre-edit) values for rows 11-13 in the affected columns.
$orig = @{}
foreach ($row in 11..13) {
    $orig[$row] = @{}
    foreach ($col in $cols) {
        $orig[$row][$col] = $ws.Range("$col$row").Value2
    }
}

# Map: target row -> source row (cyclic rotation)
$mapping = @{ 11 = 13; 12 = 11; 13 = 12 }

foreach ($targetRow in $mapping.Keys) {
    $sourceRow = $mapping[$targetRow]
    foreach ($col in $cols) {
        $ws.Range("$col$targetRow").Value = $orig[$sourceRow][$col]
    }
}
